# Generate Report for Handback
# Updates the localization status workbook to reflect that translations have
# been handed back (in sync with en-US) for the zh-cn and de-de locales:
#   - Status columns flip from "Ready for handoff" to "Handed back: in sync with en-US"
#   - New "Latest Target File" / "Latest Handback File" hyperlinked cells are populated
#   - "Latest Handback DateTime" is stamped per locale

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$targetFileName = "a121ee60-f53a-42a1-88f4-e6ee8c3ea897.md"
$targetFileUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/ba5a2b4364fc3cef8a0adec22f1c661f18676169/e2e/a121ee60-f53a-42a1-88f4-e6ee8c3ea897.md"

# ---- Overview sheet: just the status text flips ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhHandbackFile = "a121ee60-f53a-42a1-88f4-e6ee8c3ea897.81435a961352706e7c8f45955ddb767c7ab3c8b1.zh-cn.xlf"
$zhHandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9bbb176ced3a12684e130be9a8bc75db6fafa19a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a121ee60-f53a-42a1-88f4-e6ee8c3ea897.81435a961352706e7c8f45955ddb767c7ab3c8b1.zh-cn.xlf"
$zhHandbackTime = "2016-03-12 20:53:12"

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $targetFileUrl, "", "", $targetFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhHandbackUrl, "", "", $zhHandbackFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $targetFileUrl, "", "", $targetFileName)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhHandbackUrl, "", "", $zhHandbackFile)

$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Range("G3").Style = "HyperLink"

$wsZh.Range("H2").Value = $zhHandbackTime
$wsZh.Range("H3").Value = $zhHandbackTime

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$deHandbackFile = "a121ee60-f53a-42a1-88f4-e6ee8c3ea897.81435a961352706e7c8f45955ddb767c7ab3c8b1.de-de.xlf"
$deHandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/467e848ce85d67ac36afea7288c5a6ea23cd1f57/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a121ee60-f53a-42a1-88f4-e6ee8c3ea897.81435a961352706e7c8f45955ddb767c7ab3c8b1.de-de.xlf"
$deHandbackTime = "2016-03-12 20:53:18"

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $targetFileUrl, "", "", $targetFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deHandbackUrl, "", "", $deHandbackFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $targetFileUrl, "", "", $targetFileName)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deHandbackUrl, "", "", $deHandbackFile)

$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Range("G3").Style = "HyperLink"

$wsDe.Range("H2").Value = $deHandbackTime
$wsDe.Range("H3").Value = $deHandbackTime
